# Refine metadata to be an additional tab.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the existing "data" sheet ---
$data.Range("F2").Value  = "2021-10-05 14:19:08.433858"
$data.Range("F3").Value  = "2021-10-05 14:19:08.433866"
$data.Range("F4").Value  = "2021-10-05 14:19:08.433870"
$data.Range("F5").Value  = "2021-10-05 14:19:08.433873"
$data.Range("F6").Value  = "2021-10-05 14:19:08.433876"
$data.Range("F7").Value  = "2021-10-05 14:19:08.433879"
$data.Range("F8").Value  = "2021-10-05 14:19:08.433881"
$data.Range("F9").Value  = "2021-10-05 14:19:08.433885"
$data.Range("F10").Value = "2021-10-05 14:19:08.433888"
$data.Range("F11").Value = "2021-10-05 14:19:08.433891"
$data.Range("F12").Value = "2021-10-05 14:19:08.433893"

# --- Add a new "metadata" worksheet, positioned after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$meta.Name = "metadata"

# Match the bold / centered / bordered header style used on the "data" sheet.
$meta.Range("B1:G1").Font.Bold = $true
$meta.Range("B1:G1").HorizontalAlignment = -4108
$meta.Range("B1:G1").VerticalAlignment = -4160
$meta.Range("B1:G1").Borders.LineStyle = 1

$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160
$meta.Range("A2").Borders.LineStyle = 1

# Header row.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Amyloidosis"
$meta.Range("C2").Value = 502
$meta.Range("D2").Value = "'1.11"
$meta.Range("E2").Value = "2021-08-03T17:01:01.816597Z"
$meta.Range("F2").Value = "2021-10-05 14:19:08.430196"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/502/?format=json"

# Keep the originally-active "data" sheet selected/active.
$data.Activate()

